# Replace curly double quotes (U+201C / U+201D) with a straight single
# quote (U+0027) in a handful of English dialogue lines on Sheet1,
# matching commit "update on 20210731 画中人".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$leftDoubleQuote  = [char]0x201C
$rightDoubleQuote = [char]0x201D
$straightSingle   = [char]0x0027

$targetCells = @("C19", "C22", "C54", "C56", "C69", "C76")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $val = $cell.Value2
    $newVal = $val.Replace($leftDoubleQuote, $straightSingle).Replace($rightDoubleQuote, $straightSingle)
    $cell.Value = $newVal
}
